# Add two new columns, I (I0) and J (IF), to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers for the new columns, matching the style of the existing header row (H1).
$ws.Range("H1").Copy($ws.Range("I1"))
$ws.Range("H1").Copy($ws.Range("J1"))
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for I2:J29 ([I, J] per row).
$data = @(
    @(8, 8),
    @(8, 8),
    @(6, 6),
    @(5, 7),
    @(8, 8),
    @(4, 5),
    @(4, 6),
    @(8, 8),
    @(4, 6),
    @(6, 6),
    @(7, 8),
    @(7, 7),
    @(8, 8),
    @(6, 7),
    @(5, 5),
    @(8, 8),
    @(7, 8),
    @(8, 8),
    @(9, 9),
    @(4, 4),
    @(8, 9),
    @(6, 7),
    @(6, 7),
    @(7, 7),
    @(9, 9),
    @(9, 9),
    @(1, 1),
    @(8, 8)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
